$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Report" to "Sheet1"
$ws.Name = "Sheet1"

# Update membership counts (column B) for rows 3-6
$ws.Range("B3").Value = 700
$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 100
$ws.Range("B6").Value = 1000
